$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B289/D289
$ws.Cells.Item(289, 2).Value = 32957600000
$ws.Cells.Item(289, 4).Value = 46550282485.87571

# Update B322/D322
$ws.Cells.Item(322, 2).Value = 36383400000
$ws.Cells.Item(322, 4).Value = 51403503814.6369

# Update B327/D327
$ws.Cells.Item(327, 2).Value = 37472200000
$ws.Cells.Item(327, 4).Value = 52941791466.51597

# Update B328/D328
$ws.Cells.Item(328, 2).Value = 37492700000
$ws.Cells.Item(328, 4).Value = 52970754450.40972

# Add new rows 352/353, copying formatting style from row 351
$ws.Range("A351:D351").Copy()
$ws.Range("A352:D353").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(352, 1).Value = 44986
$ws.Cells.Item(352, 2).Value = 41710700000
$ws.Cells.Item(352, 3).Value = 1.409443269908386
$ws.Cells.Item(352, 4).Value = 58788865398.16772

$ws.Cells.Item(353, 1).Value = 45017
$ws.Cells.Item(353, 2).Value = 41839100000
$ws.Cells.Item(353, 3).Value = 1.410835214446953
$ws.Cells.Item(353, 4).Value = 59028075620.7675

Write-Host "Done"
